$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was added at the top of the data set (row 2),
# pushing every existing record down by one row.
$ws.Rows(2).Insert()

# The freshly inserted row comes in with the formatting copied from the
# row above (the bold header). Reset it to the plain/default style used
# by every other data row, then restore the date number format on the
# "Fecha" column (D), matching the rest of the table.
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = "2021-11-15"
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 100112026
$ws.Range("G2").Value = "Haba"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 155
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 9000
$ws.Range("N2").Value = "$/saco 25 kilos"
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("P2").Value = 360
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
